$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(421).Insert()

$ws.Cells.Item(421, 1).Value = 11
$ws.Cells.Item(421, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(421, 3).Value = "Bíobío"
$ws.Cells.Item(421, 4).Value = 45212
$ws.Cells.Item(421, 5).Value = 8
$ws.Cells.Item(421, 6).Value = 100112045
$ws.Cells.Item(421, 7).Value = "Zapallo"
$ws.Cells.Item(421, 8).Value = "Camote"
$ws.Cells.Item(421, 9).Value = "1a nueva(o)"
$ws.Cells.Item(421, 10).Value = 500
$ws.Cells.Item(421, 11).Value = 1000
$ws.Cells.Item(421, 12).Value = 1000
$ws.Cells.Item(421, 13).Value = 1000
$ws.Cells.Item(421, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(421, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(421, 16).Value = 1000
$ws.Cells.Item(421, 17).Value = 1
$ws.Cells.Item(421, 18).Value = "Hortaliza"

$ws.Cells.Item(421, 4).NumberFormat = $ws.Cells.Item(422, 4).NumberFormat
